$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header "SELECTED_MINIMISATION" in column R (R1),
# matching the style used by the other headers (e.g. Q1).
$ws.Range("R1").Value = "SELECTED_MINIMISATION"
$ws.Range("R1").HorizontalAlignment = $ws.Range("Q1").HorizontalAlignment

# Move the active selection to R2, as recorded in the saved view state.
$ws.Range("R2").Select()
